$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))
$ws.Cells.Item($row, 1).Value = 45986

$ws.Cells.Item($row, 2).Value = 2025
$ws.Cells.Item($row, 3).Value = 0.2298740481777584
$ws.Cells.Item($row, 4).Value = 2026
$ws.Cells.Item($row, 5).Value = -0.05255865067609333
